# Extended unified config process to EAs and OPs
# - Rename the two worksheets to reflect their new purpose
# - Make the (renamed) second sheet ("OPs_cat") the active/selected tab,
#   matching the workbook's saved view state in the target file.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "OPs_new"
$ws2.Name = "OPs_cat"

# Switch the active tab from sheet 1 to sheet 2 (OPs_cat), which updates
# workbookView/activeTab and moves sheetView/tabSelected accordingly.
$ws2.Activate()
